# Employee_Data.xlsx update:
#  - strip the stray literal quote characters that were embedded in the
#    branch-name strings (column C, rows 2-16)
#  - switch the hire_date column (F) from the long custom date format to a
#    short m/d/yyyy date format
#  - append 15 new employee rows (17-31) for a second branch (2117272),
#    mirroring the layout/format of the existing rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Clean up the branch-name text in column C (remove the extra quote marks)
# ---------------------------------------------------------------------------
$branchNames = @{
    2  = "Bank Sinatra"
    3  = "Bank Langella"
    4  = "Anne Bank"
    5  = "Bank and File"
    6  = "Draw a Bank"
    7  = "Point-Bank"
    8  = "Walk the Bank"
    9  = "Bank Zappa"
    10 = "Bank Ocean"
    11 = "Bank Ocean"
    12 = "Bank Sinatra"
    13 = "Purple Bank"
    14 = "Bank You"
    15 = "Swiss Bank"
    16 = "Military Bank"
}
foreach ($r in $branchNames.Keys) {
    $ws.Range("C$r").Value = $branchNames[$r]
}

# ---------------------------------------------------------------------------
# 2. Re-format the hire_date column to a short date (creates/uses the new
#    numFmtId 14 style) and propagate it down F2:F16 via a format-only paste
#    so every cell shares the same style index.
# ---------------------------------------------------------------------------
$ws.Range("F2").NumberFormat = "mm-dd-yy"
$ws.Range("F2").Copy()
$ws.Range("F3:F16").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Append the new rows (17-31). First clone row 16's formatting down so the
#    new rows pick up the right (shared) styles for E/F, then fill in values.
# ---------------------------------------------------------------------------
$ws.Range("A16:F16").Copy()
$ws.Range("A17:F31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newRows = @(
    @{ Row=17; B=45; C="Money Man";        D="Branch Manager";                 E=2117272; F="2017-08-01" },
    @{ Row=18; B=92; C="Money Man Two";    D="Assistant Branch Manager";       E=2117272; F="2017-08-01" },
    @{ Row=19; B=34; C="Mo Money";         D="Loan Officer";                   E=2117272; F="2017-08-13" },
    @{ Row=20; B=54; C="Money Honey";      D="Teller";                         E=2117272; F="2017-08-13" },
    @{ Row=21; B=65; C="Bands of Money";   D="Teller";                         E=2117272; F="2017-09-02" },
    @{ Row=22; B=85; C="Fake Money";       D="Teller";                         E=2117272; F="2017-10-16" },
    @{ Row=23; B=23; C="No Money";         D="Teller";                         E=2117272; F="2018-08-07" },
    @{ Row=24; B=42; C="Money Money";      D="Teller";                         E=2117272; F="2017-08-08" },
    @{ Row=25; B=55; C="Money";            D="Personal Banker";                E=2117272; F="2019-08-09" },
    @{ Row=26; B=43; C="Less Money";       D="Personal Banker";                E=2117272; F="2017-08-10" },
    @{ Row=27; B=22; C="Money Less";       D="Personal Banker";                E=2117272; F="2017-08-11" },
    @{ Row=28; B=21; C="Who Has Money";    D="Customer Service Representative"; E=2117272; F="2020-01-12" },
    @{ Row=29; B=19; C="Out Of Money";     D="Financial Advisor";              E=2117272; F="2017-08-13" },
    @{ Row=30; B=23; C="My Money";         D="Financial Advisor";              E=2117272; F="2017-08-14" },
    @{ Row=31; B=40; C="Our Money";        D="Financial Advisor";              E=2117272; F="2017-08-15" }
)

foreach ($row in $newRows) {
    $r = $row.Row
    $prev = $r - 1
    $ws.Range("A$r").Formula = "=A$prev+1"
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    if ($r -eq 17) {
        $ws.Range("D$r").Value = "'" + $row.D
    } else {
        $ws.Range("D$r").Value = $row.D
    }
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = (Get-Date -Date $row.F)
}

$ws.Range("G10").Select()
